# Auto-generated: applies scheduled market-data refresh to Omega_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 11537
$ws.Range("I76").Value = 11816.7
$ws.Range("J76").Value = 10915.444
$ws.Range("K76").Value = 11816.7
$ws.Range("L76").Value = 10915.444
$ws.Range("M76").Value = -11501.7
$ws.Range("N76").Value = -11545.444

$ws.Range("H79").Value = 11537
$ws.Range("I79").Value = 11816.7
$ws.Range("J79").Value = 10915.444
$ws.Range("K79").Value = 11816.7
$ws.Range("L79").Value = 10915.444
$ws.Range("M79").Value = -10724.7
$ws.Range("N79").Value = -13099.444

$ws.Range("H96").Value = 2294.5
$ws.Range("I96").Value = 589
$ws.Range("K96").Value = 1767
$ws.Range("M96").Value = -394

$ws.Range("H98").Value = 1605.125
$ws.Range("I98").Value = 1605.125
$ws.Range("K98").Value = 1605.125
$ws.Range("M98").Value = -107.125

$ws.Range("H103").Value = 764.4
$ws.Range("I103").Value = 394.75
$ws.Range("J103").Value = 1010.8333
$ws.Range("K103").Value = 1184.25
$ws.Range("L103").Value = 3032.4999
$ws.Range("M103").Value = -598.25
$ws.Range("N103").Value = -4204.4999

$ws.Range("H111").Value = 11406.6
$ws.Range("I111").Value = 10415.625
$ws.Range("J111").Value = 12241.105
$ws.Range("K111").Value = 31246.875
$ws.Range("L111").Value = 36723.315
$ws.Range("M111").Value = -28179.875
$ws.Range("N111").Value = -42857.315

$ws.Range("H122").Value = 1605.125
$ws.Range("I122").Value = 1605.125
$ws.Range("K122").Value = 4815.375
$ws.Range("M122").Value = -2365.375

$ws.Range("H140").Value = 114495
$ws.Range("J140").Value = 114495
$ws.Range("L140").Value = 114495
$ws.Range("N140").Value = -124855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 5333
$ws.Range("I19").Value = 5333
$ws.Range("K19").Value = 5333
$ws.Range("M19").Value = -5104

$ws.Range("H45").Value = 10798.962
$ws.Range("I45").Value = 13538.277
$ws.Range("K45").Value = 13538.277
$ws.Range("M45").Value = -13161.277

$ws.Range("H74").Value = 2625.1765
$ws.Range("I74").Value = 2662.5334
$ws.Range("K74").Value = 2662.5334
$ws.Range("M74").Value = -1788.5334

$ws.Range("H77").Value = 2625.1765
$ws.Range("I77").Value = 2662.5334
$ws.Range("K77").Value = 13312.667
$ws.Range("M77").Value = -8944.666999999999

$ws.Range("H97").Value = 1965.75
$ws.Range("I97").Value = 1844.4546
$ws.Range("J97").Value = 3300
$ws.Range("K97").Value = 1844.4546
$ws.Range("L97").Value = 3300
$ws.Range("M97").Value = -1348.4546
$ws.Range("N97").Value = -4292

$ws.Range("N102").ClearContents()
$ws.Range("H102").Value = 5253.1665
$ws.Range("I102").Value = 5253.1665
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5253.1665
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3631.1665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N59").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0

$ws.Range("H105").Value = 4162.364
$ws.Range("I105").Value = 4320.6665
$ws.Range("K105").Value = 4320.6665
$ws.Range("M105").Value = -2573.6665

$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0

$ws.Range("H138").Value = 89214.375
$ws.Range("J138").Value = 89214.375
$ws.Range("L138").Value = 89214.375
$ws.Range("N138").Value = -99494.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3318.5789
$ws.Range("I31").Value = 4143.35
$ws.Range("J31").Value = 2402.1667
$ws.Range("K31").Value = 4143.35
$ws.Range("L31").Value = 2402.1667
$ws.Range("M31").Value = -3848.35
$ws.Range("N31").Value = -2992.1667

$ws.Range("H34").Value = 3318.5789
$ws.Range("I34").Value = 4143.35
$ws.Range("J34").Value = 2402.1667
$ws.Range("K34").Value = 4143.35
$ws.Range("L34").Value = 2402.1667
$ws.Range("M34").Value = -3941.35
$ws.Range("N34").Value = -2806.1667

$ws.Range("H134").Value = 8461.619000000001
$ws.Range("I134").Value = 7332.3887
$ws.Range("K134").Value = 21997.1661
$ws.Range("M134").Value = -19462.1661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2765.625
$ws.Range("I3").Value = 2364.1738
$ws.Range("K3").Value = 7092.5214
$ws.Range("M3").Value = -6980.5214

$ws.Range("H107").Value = 1087.625
$ws.Range("J107").Value = 1902
$ws.Range("L107").Value = 5706
$ws.Range("N107").Value = -9546

$ws.Range("H114").Value = 1873.5714
$ws.Range("J114").Value = 2257.4
$ws.Range("L114").Value = 6772.200000000001
$ws.Range("N114").Value = -13280.2

$ws.Range("H140").Value = 1700.25
$ws.Range("I140").Value = 1428
$ws.Range("J140").Value = 3969
$ws.Range("K140").Value = 4284
$ws.Range("L140").Value = 11907
$ws.Range("M140").Value = 896
$ws.Range("N140").Value = -22267

$ws.Range("H141").Value = 3712
$ws.Range("I141").Value = 1616.3334
$ws.Range("K141").Value = 4849.0002
$ws.Range("M141").Value = 330.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4376098.5
$ws.Range("I11").Value = 213256.7
$ws.Range("J11").Value = 7349557.5
$ws.Range("K11").Value = 213256.7
$ws.Range("L11").Value = 7349557.5
$ws.Range("M11").Value = -213117.7
$ws.Range("N11").Value = -7349835.5

$ws.Range("H18").Value = 73335
$ws.Range("I18").Value = 20005
$ws.Range("K18").Value = 20005
$ws.Range("M18").Value = -19712

$ws.Range("H80").Value = 1902
$ws.Range("I80").Value = 1899
$ws.Range("J80").Value = 1905
$ws.Range("K80").Value = 1899
$ws.Range("L80").Value = 1905
$ws.Range("M80").Value = -901
$ws.Range("N80").Value = -3901

$ws.Range("H83").Value = 1902
$ws.Range("I83").Value = 1899
$ws.Range("J83").Value = 1905
$ws.Range("K83").Value = 9495
$ws.Range("L83").Value = 9525
$ws.Range("M83").Value = -4503
$ws.Range("N83").Value = -19509

$ws.Range("H122").Value = 3737.25
$ws.Range("I122").Value = 3475
$ws.Range("K122").Value = 10425
$ws.Range("M122").Value = -7975

$ws.Range("H126").Value = 6544.643
$ws.Range("I126").Value = 5321.7144
$ws.Range("K126").Value = 15965.1432
$ws.Range("M126").Value = -13495.1432

$ws.Range("H132").Value = 4615.0713
$ws.Range("I132").Value = 5143.7896
$ws.Range("K132").Value = 15431.3688
$ws.Range("M132").Value = -12901.3688

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1308
$ws.Range("I22").Value = 1409.8
$ws.Range("K22").Value = 1409.8
$ws.Range("M22").Value = -1114.8

$ws.Range("H23").Value = 1825.25
$ws.Range("I23").Value = 1825.25
$ws.Range("K23").Value = 1825.25
$ws.Range("M23").Value = -1595.25

$ws.Range("H27").Value = 1308
$ws.Range("I27").Value = 1409.8
$ws.Range("K27").Value = 1409.8
$ws.Range("M27").Value = -1302.8

$ws.Range("H114").Value = 88998.664
$ws.Range("J114").Value = 88998.664
$ws.Range("L114").Value = 88998.664
$ws.Range("N114").Value = -97676.664

$ws.Range("H122").Value = 8468.429
$ws.Range("I122").Value = 9296.5
$ws.Range("K122").Value = 27889.5
$ws.Range("M122").Value = -25439.5

$ws.Range("H132").Value = 2501.5
$ws.Range("I132").Value = 2501.5
$ws.Range("K132").Value = 7504.5
$ws.Range("M132").Value = -4974.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3301
$ws.Range("I81").Value = 600
$ws.Range("J81").Value = 6002
$ws.Range("K81").Value = 1200
$ws.Range("L81").Value = 12004
$ws.Range("M81").Value = -139
$ws.Range("N81").Value = -14126

$ws.Range("H84").Value = 3301
$ws.Range("I84").Value = 600
$ws.Range("J84").Value = 6002
$ws.Range("K84").Value = 6000
$ws.Range("L84").Value = 60020
$ws.Range("M84").Value = -696
$ws.Range("N84").Value = -70628

$ws.Range("H122").Value = 3345.1
$ws.Range("I122").Value = 3400.6
$ws.Range("J122").Value = 3289.6
$ws.Range("K122").Value = 10201.8
$ws.Range("L122").Value = 9868.799999999999
$ws.Range("M122").Value = -7751.799999999999
$ws.Range("N122").Value = -14768.8

$ws.Range("H132").Value = 3001.3794
$ws.Range("I132").Value = 3239.6
$ws.Range("J132").Value = 1909.6
$ws.Range("K132").Value = 9718.799999999999
$ws.Range("M132").Value = -7188.799999999999

